$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45972, 9101.36445167059, 8651.38430844298, 13252, 4897.30130111058, 12.3619003980647),
    @(3, 45973, 9000.66721514185, 8453.56254005483, 9004, 5054.93789734953, 187.687518225182),
    @(4, 45974, 8923.80339350939, 8237.1621090861, 9004, 4987.80619106812, 175.873679173093),
    @(5, 45975, 8819.86633177528, 7296.59938002071, 9004, 4904.60735504893, 133.216947294568),
    @(6, 45976, 2982.08271944182, 4749.82024908422, 9004, 4265.42144002794, 0.46840371300641),
    @(7, 45977, 3412.56435984501, 4787.39839778532, 9004, 4921.19116164753, 29.3578983097023),
    @(8, 45978, 10116.508538002, 8041.33721748613, 9004, 5829.39572450412, 202.780539249594),
    @(9, 45979, 10116.508538002, 8791.9626816311, 9004, 5829.39572450412, 234.056600255634),
    @(10, 45980, 10116.508538002, 9289.75704198764, 9004, 5829.39572450412, 254.798031937157),
    @(11, 45981, 10116.508538002, 9504.95369375453, 9004, 5829.39572450412, 263.76455909411),
    @(12, 45982, 10116.508538002, 9042.0662558458, 9004, 5829.39572450412, 244.47758251458),
    @(13, 45983, 4000.30754528634, 6539.66869020191, 9004, 5477.55666137121, 125.551056315547),
    @(14, 45984, 3853.28814954241, 6194.81609841946, 9004, 5469.55316195516, 110.848719182276),
    @(15, 45985, 10980.8696790314, 9493.61654658831, 9004, 6409.03351963316, 287.443752759228)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
